$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 600.6316
$ws.Range("I19").Value = 585.7143
$ws.Range("J19").Value = 609.3333
$ws.Range("K19").Value = 585.7143
$ws.Range("L19").Value = 609.3333
$ws.Range("M19").Value = -410.7143
$ws.Range("N19").Value = -959.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 860484.6
$ws.Range("I98").Value = 1015936.4
$ws.Range("J98").Value = 5500
$ws.Range("K98").Value = 1015936.4
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = -1014438.4
$ws.Range("N98").Value = -8496

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 654148.4
$ws.Range("I107").Value = 794198.5
$ws.Range("K107").Value = 794198.5
$ws.Range("M107").Value = -792278.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 860484.6
$ws.Range("I122").Value = 1015936.4
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 3047809.2
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -3045359.2
$ws.Range("N122").Value = -21400

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 27779100
$ws.Range("I137").Value = 41667616
$ws.Range("J137").Value = 2069
$ws.Range("K137").Value = 125002848
$ws.Range("L137").Value = 6207
$ws.Range("M137").Value = -125000298
$ws.Range("N137").Value = -11307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 42126.64
$ws.Range("I2").Value = 65110.312
$ws.Range("J2").Value = 1266.7778
$ws.Range("K2").Value = 65110.312
$ws.Range("L2").Value = 1266.7778
$ws.Range("M2").Value = -64997.312
$ws.Range("N2").Value = -1492.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4717.6665
$ws.Range("I88").Value = 1600
$ws.Range("J88").Value = 6276.5
$ws.Range("K88").Value = 1600
$ws.Range("L88").Value = 6276.5
$ws.Range("M88").Value = -1194
$ws.Range("N88").Value = -7088.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 4717.6665
$ws.Range("I91").Value = 1600
$ws.Range("J91").Value = 6276.5
$ws.Range("K91").Value = 1600
$ws.Range("L91").Value = 6276.5
$ws.Range("M91").Value = -196
$ws.Range("N91").Value = -9084.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 272877
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 272877
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 272877
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -275651

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1445.4166
$ws.Range("I110").Value = 693.0769
$ws.Range("K110").Value = 693.0769
$ws.Range("M110").Value = 1351.9231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 42126.64
$ws.Range("I116").Value = 65110.312
$ws.Range("J116").Value = 1266.7778
$ws.Range("K116").Value = 65110.312
$ws.Range("L116").Value = 1266.7778
$ws.Range("M116").Value = -62816.312
$ws.Range("N116").Value = -5854.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2962.4
$ws.Range("I122").Value = 3078
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 9234
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -6784
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 42126.64
$ws.Range("I3").Value = 65110.312
$ws.Range("J3").Value = 1266.7778
$ws.Range("K3").Value = 65110.312
$ws.Range("L3").Value = 1266.7778
$ws.Range("M3").Value = -64996.312
$ws.Range("N3").Value = -1494.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 900
$ws.Range("J64").Value = 900
$ws.Range("L64").Value = 900
$ws.Range("N64").Value = -1350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 900
$ws.Range("J67").Value = 900
$ws.Range("L67").Value = 900
$ws.Range("N67").Value = -2460

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1189.7059
$ws.Range("I107").Value = 878.1818
$ws.Range("K107").Value = 878.1818
$ws.Range("M107").Value = 1041.8182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3850.8635
$ws.Range("I134").Value = 2542.25
$ws.Range("J134").Value = 5421.2
$ws.Range("K134").Value = 7626.75
$ws.Range("L134").Value = 16263.6
$ws.Range("M134").Value = -5091.75
$ws.Range("N134").Value = -21333.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2749.9167
$ws.Range("I5").Value = 540.4286
$ws.Range("J5").Value = 5843.2
$ws.Range("K5").Value = 540.4286
$ws.Range("L5").Value = 5843.2
$ws.Range("M5").Value = -428.4286
$ws.Range("N5").Value = -6067.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 64
$ws.Range("I7").Value = 67.5
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 67.5
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 45.5
$ws.Range("N7").Value = -276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 7859.2856
$ws.Range("I8").Value = 5500
$ws.Range("J8").Value = 8803
$ws.Range("K8").Value = 5500
$ws.Range("L8").Value = 8803
$ws.Range("M8").Value = -5360
$ws.Range("N8").Value = -9083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 14006.75
$ws.Range("J14").Value = 14006.75
$ws.Range("L14").Value = 14006.75
$ws.Range("N14").Value = -14346.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1259.2222
$ws.Range("I16").Value = 1174.9286
$ws.Range("J16").Value = 1554.25
$ws.Range("K16").Value = 1174.9286
$ws.Range("L16").Value = 1554.25
$ws.Range("M16").Value = -887.9286
$ws.Range("N16").Value = -2128.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 958.75
$ws.Range("I19").Value = 381.42856
$ws.Range("K19").Value = 381.42856
$ws.Range("M19").Value = -211.42856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 34801.6
$ws.Range("I23").Value = 31002.25
$ws.Range("J23").Value = 49999
$ws.Range("K23").Value = 31002.25
$ws.Range("L23").Value = 49999
$ws.Range("M23").Value = -30762.25
$ws.Range("N23").Value = -50479

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 958.75
$ws.Range("I24").Value = 381.42856
$ws.Range("K24").Value = 381.42856
$ws.Range("M24").Value = -211.42856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 5248.3335
$ws.Range("I25").Value = 4997.778
$ws.Range("J25").Value = 6000
$ws.Range("K25").Value = 4997.778
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = -4823.778
$ws.Range("N25").Value = -6348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 34801.6
$ws.Range("I27").Value = 31002.25
$ws.Range("J27").Value = 49999
$ws.Range("K27").Value = 31002.25
$ws.Range("L27").Value = 49999
$ws.Range("M27").Value = -30810.25
$ws.Range("N27").Value = -50383

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2409.5217
$ws.Range("I31").Value = 1591.2142
$ws.Range("J31").Value = 3682.4443
$ws.Range("K31").Value = 1591.2142
$ws.Range("L31").Value = 3682.4443
$ws.Range("M31").Value = -1296.2142
$ws.Range("N31").Value = -4272.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2409.5217
$ws.Range("I34").Value = 1591.2142
$ws.Range("J34").Value = 3682.4443
$ws.Range("K34").Value = 1591.2142
$ws.Range("L34").Value = 3682.4443
$ws.Range("M34").Value = -1389.2142
$ws.Range("N34").Value = -4086.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2291.0417
$ws.Range("I58").Value = 1151.6875
$ws.Range("J58").Value = 4569.75
$ws.Range("K58").Value = 1151.6875
$ws.Range("L58").Value = 4569.75
$ws.Range("M58").Value = -948.6875
$ws.Range("N58").Value = -4975.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 320.42856
$ws.Range("I107").Value = 218.6923
$ws.Range("J107").Value = 485.75
$ws.Range("K107").Value = 218.6923
$ws.Range("L107").Value = 485.75
$ws.Range("M107").Value = 1701.3077
$ws.Range("N107").Value = -4325.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1259.2222
$ws.Range("I113").Value = 1174.9286
$ws.Range("J113").Value = 1554.25
$ws.Range("K113").Value = 1174.9286
$ws.Range("L113").Value = 1554.25
$ws.Range("M113").Value = 995.0714
$ws.Range("N113").Value = -5894.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2291.0417
$ws.Range("I136").Value = 1151.6875
$ws.Range("J136").Value = 4569.75
$ws.Range("K136").Value = 3455.0625
$ws.Range("L136").Value = 13709.25
$ws.Range("M136").Value = -905.0625
$ws.Range("N136").Value = -18809.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1750
$ws.Range("J57").Value = 4000
$ws.Range("L57").Value = 12000
$ws.Range("N57").Value = -13118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2144.4167
$ws.Range("I114").Value = 1021
$ws.Range("J114").Value = 2706.125
$ws.Range("K114").Value = 3063
$ws.Range("L114").Value = 8118.375
$ws.Range("M114").Value = 191
$ws.Range("N114").Value = -14626.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 7983.875
$ws.Range("I99").Value = 6267.2856
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 6267.2856
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = -4021.2856
$ws.Range("N99").Value = -24492

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2223.647
$ws.Range("I113").Value = 2190.2
$ws.Range("J113").Value = 2271.4285
$ws.Range("K113").Value = 2190.2
$ws.Range("L113").Value = 2271.4285
$ws.Range("M113").Value = -20.19999999999982
$ws.Range("N113").Value = -6611.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1650
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 1725
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 1725
$ws.Range("M61").Value = -1298
$ws.Range("N61").Value = -2129

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1650
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1725
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1725
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6065

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 341.16666
$ws.Range("I113").Value = 258
$ws.Range("J113").Value = 507.5
$ws.Range("K113").Value = 774
$ws.Range("L113").Value = 1522.5
$ws.Range("M113").Value = 1396
$ws.Range("N113").Value = -5862.5
